$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate each card's multi-row details into a single row per card,
# formatted as a Python-tuple-like string: ('Name', ['field1', 'field2', ...])

$ws.Range("A2").Value = "('Basandra, Battle Seraph', ['{3}{R}{W}', 'Legendary Creature — Angel', 'Flying', 'Players can" + [char]8217 + "t cast spells during combat.', '{R}: Target creature attacks this turn if able.', '4/4'])"

$ws.Range("A3").Value = "('Edric, Spymaster of Trest', ['{1}{G}{U}', 'Legendary Creature — Elf Rogue', 'Whenever a creature deals combat damage to one of your opponents, its controller may draw a card.', '2/2'])"

$ws.Range("A4").Value = "('Nin, the Pain Artist', ['{U}{R}', 'Legendary Creature — Vedalken Wizard', '{X}{U}{R}, {T}: Nin, the Pain Artist deals X damage to target creature. That creature" + [char]8217 + "s controller draws X cards.', '1/1'])"

$ws.Range("A5").Value = "('Skullbriar, the Walking Grave', ['{B}{G}', 'Legendary Creature — Zombie Elemental', 'Haste', 'Whenever Skullbriar, the Walking Grave deals combat damage to a player, put a +1/+1 counter on it.', 'Counters remain on Skullbriar as it moves to any zone other than a player" + [char]8217 + "s hand or library.', '1/1'])"

$ws.Range("A6").Value = "('Vish Kal, Blood Arbiter', ['{4}{W}{B}{B}', 'Legendary Creature — Vampire', 'Flying, lifelink', 'Sacrifice a creature: Put X +1/+1 counters on Vish Kal, Blood Arbiter, where X is the sacrificed creature" + [char]8217 + "s power.', 'Remove all +1/+1 counters from Vish Kal: Target creature gets -1/-1 until end of turn for each +1/+1 counter removed this way.', '5/5'])"

# Remove the now-unused rows 7:32 (previously held the split-out details)
$rows = $ws.Range("A7:A32").EntireRow
$rows.Delete()
